$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert a new "2022-Q1" worksheet (fund holdings detail) right
#    before the "总计" (totals) sheet. We duplicate an existing
#    quarterly sheet ("2021-Q3") so the new sheet starts out with the
#    same layout/formatting, then we overwrite its contents.
# ---------------------------------------------------------------------
$src = $wb.Worksheets.Item("2021-Q3")
$total = $wb.Worksheets.Item("总计")
$src.Copy($total)
$q1 = $wb.ActiveSheet
$q1.Name = "2022-Q1"

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$q1.Range("A2").Value = 0
$q1.Range("B2").Value = "'000059"
$q1.Range("C2").Value = "国联安中证医药100指数A"
$q1.Range("D2").Value = "'2.20"
$q1.Range("E2").Value = "'91.27"
$q1.Range("F2").Value = "'1.28"
$q1.Range("G2").Value = "'0.0282"
$q1.Range("H2").Value = 9

$q1.Range("A3").Value = 1
$q1.Range("B3").Value = "'006569"
$q1.Range("C3").Value = "国联安中证医药100指数C"
$q1.Range("D3").Value = "'0.34"
$q1.Range("E3").Value = "'91.27"
$q1.Range("F3").Value = "'1.28"
$q1.Range("G3").Value = "'0.0044"
$q1.Range("H3").Value = 9

# ---------------------------------------------------------------------
# 2) Add the corresponding "2022-Q1" row to the "总计" summary sheet,
#    above the existing rows (which shift down by one).
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("总计")
$ws.Rows.Item(2).Insert()

$ws.Range("B2:D2").ClearFormats()

$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "2022-Q1"
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 0.03

$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4
